# Apply the "season record" columns (Wins / Losses / Ties) to the BAL_2022
# worksheet. Header row 1 gets three new bold/bordered header cells in
# AD1:AF1, and every data row (2-60) gets the team's season record
# (83 wins, 79 losses, 0 ties) in AD:AF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells: copy the existing header formatting (bold font, thin
# border, centered alignment) from column A's header onto the new header
# cells, then set their text. ---
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data cells: every player row shares the same team record. ---
$wins = 83
$losses = 79
$ties = 0

$lastRow = 60
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
